$d = $word.ActiveDocument

$pairs = @(
    @{Old = "147×6=882";   New = "582×5=2910"},
    @{Old = "213×5=1065";  New = "705×5=3525"},
    @{Old = "130×8=1040";  New = "440×6=2640"},
    @{Old = "209×7=1463";  New = "141×8=1128"},
    @{Old = "250×7=1750";  New = "384×4=1536"},
    @{Old = "648×8=5184";  New = "965×3=2895"},
    @{Old = "221×4=884";   New = "732×3=2196"},
    @{Old = "910×5=4550";  New = "856×3=2568"},
    @{Old = "761×4=3044";  New = "694×7=4858"},
    @{Old = "743×6=4458";  New = "209×3=627"},
    @{Old = "887×5=4435";  New = "612×4=2448"},
    @{Old = "983×2=1966";  New = "875×5=4375"},
    @{Old = "341×7=2387";  New = "625×5=3125"},
    @{Old = "111×9=999";   New = "261×5=1305"},
    @{Old = "922×7=6454";  New = "246×2=492"},
    @{Old = "911×2=1822";  New = "353×2=706"},
    @{Old = "908×2=1816";  New = "313×9=2817"},
    @{Old = "169×5=845";   New = "574×3=1722"},
    @{Old = "919×7=6433";  New = "452×5=2260"},
    @{Old = "311×2=622";   New = "297×7=2079"},
    @{Old = "249×7=1743";  New = "796×4=3184"},
    @{Old = "864×5=4320";  New = "927×5=4635"},
    @{Old = "751×6=4506";  New = "364×8=2912"},
    @{Old = "684×8=5472";  New = "991×9=8919"},
    @{Old = "114×2=228";   New = "260×8=2080"}
)

foreach ($pair in $pairs) {
    $range = $d.Content
    $range.Find.Execute($pair.Old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $pair.New, 2)
}
